# Append a new "gold price" entry row to the bottom of the data table.
#
# The sheet has one row per day: column A holds a date-like label
# (e.g. "11-11-2025") stored as literal text, column B holds the long
# price description, also stored as literal text. Column A is General
# formatted, so assigning a date-shaped string straight to .Value makes
# Excel "smart"-convert it into a real date serial (and mint a brand new
# number-format style in the process). To keep the new row byte-for-byte
# consistent with its neighbours (same shared-string storage, same
# column-default style, no new style entries) we instead build the text
# via a literal formula ("=""12-11-2025""") and then collapse that
# formula down to its plain cached value with a copy / paste-values,
# exactly like pressing F9 then "Paste Values" in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Range("A56").Row + 1

$dateCell  = $ws.Cells.Item($newRow, 1)
$priceCell = $ws.Cells.Item($newRow, 2)

$newDate = "12-11-2025"
$newPrice = "The price of gold in India today is ₹12,551 per gram for 24 karat gold, ₹11,505 per gram for 22 karat gold and ₹9,413 per gram for 18 karat gold (also called 999 gold)."

# Write the date-like text as a quoted-string formula so it is never
# reinterpreted as a date, then flatten the formula to a plain cached
# text value in place (preserves the column's default General style).
$dateCell.Formula = '="' + $newDate + '"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

# Plain text, not date-shaped, so a direct value assignment is safe.
$priceCell.Value = $newPrice
